$d = $word.ActiveDocument

# 1. Version number bump: 5 -> 6
$d.Content.Find.Execute("Version 5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version 6", 2)

# 2. Date/time field text update
$d.Content.Find.Execute("5/19/22 2:15 PM", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5/31/22 7:19 AM", 2)

# 3. Typo correction: "ranked by the risk SME" -> "risked by the risk SME"
$d.Content.Find.Execute("ranked by the security SME and then ranked by the risk SME", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ranked by the security SME and then risked by the risk SME", 2)
